$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reverse the Olympiad rows so the table runs chronologically
#     (Salt Lake City 2002 -> Beijing 2022) instead of reverse-chronologically,
#     ready for plotting the dynamics of gold-medal achievements over time. ---

$data = @(
    @("Salt Lake City (2002)", 25, 36, 34, 16, 8, 17, 11),
    @("Turin (2006)",          19, 29, 24, 24, 9, 23, 14),
    @("Vancouver (2010)",      23, 29, 36, 25, 8, 16, 9),
    @("Sochi (2014)",          25, 17, 26, 23, 24, 17, 11),
    @("Pyeongchang (2018)",    36, 29, 21, 24, 20, 12, 13),
    @("Beijing (2022)",        34, 24, 21, 24, 17, 16, 15)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    for ($c = 1; $c -le 7; $c++) {
        $ws.Cells.Item($r, $c + 1).Value2 = $row[$c]
    }
    $r = $r + 1
}

# --- Make the header row bold (it was regular weight before) ---
$ws.Range("B1:G1").Font.Bold = $true
$ws.Range("H1").Font.Bold = $true

# --- Page setup tweak that came along with the fix ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Restore the active selection left after the edit ---
$ws.Range("C9").Select() | Out-Null
